$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend "ui/assets/" to the icon paths in column K for rows 6-14
for ($row = 6; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 11)  # Column K is the 11th column
    $current = $cell.Value()
    if ($current -and -not $current.StartsWith("ui/assets/")) {
        $cell.Value = "ui/assets/" + $current
    }
}
